$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Berenjena needs to be inserted as the new
# row 94 (Region de Arica y Parinacota's Huasco-province entry for
# 2021-10-22), pushing the existing rows 94-111 down to 95-112.
$ws.Rows.Item(94).Insert()

$ws.Cells.Item(94, 1).Value = 6
$ws.Cells.Item(94, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(94, 3).Value = "Metropolitana"
$ws.Cells.Item(94, 4).Value = 44491
$ws.Cells.Item(94, 5).Value = 13
$ws.Cells.Item(94, 6).Value = 100112001
$ws.Cells.Item(94, 7).Value = "Berenjena"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 220
$ws.Cells.Item(94, 11).Value = 8000
$ws.Cells.Item(94, 12).Value = 9000
$ws.Cells.Item(94, 13).Value = 8545
$ws.Cells.Item(94, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(94, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(94, 16).Value = 142
$ws.Cells.Item(94, 17).Value = 60
$ws.Cells.Item(94, 18).Value = "Hortaliza"
